# Update profit files after running on 2025-10-27
# Append the newest day's row (date + profit) to the bottom of the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 71

# Column A holds the date as plain text (matching the existing rows), so
# force a text number format first to stop Excel from auto-converting the
# "MM/DD/YYYY" string into a date serial value, then restore the default
# "Normal" style so the new cell isn't left with an explicit style index.
$dateCell = $ws.Cells.Item($newRow, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "10/27/2025"
$dateCell.Style = "Normal"

# Column B holds the numeric profit value for that date.
$ws.Cells.Item($newRow, 2).Value = 11677.36
